$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header F1 from "cuenta" to "Tipo de Producto"
$ws.Range("F1").Value = "Tipo de Producto"

# Resize column F to fit the new, longer header text (bestFit-style width)
$ws.Columns.Item(6).ColumnWidth = 14.25

# Move the active selection to G14 as recorded in the saved view state
$ws.Range("G14").Select() | Out-Null
